$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename columns to snake_case machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize "de/del/de la/de los" connector casing to "De/Del/De La/De Los"
$ws.Range("B6").Value = "Amatenango De La Frontera"
$ws.Range("B7").Value = "Comitán De Domínguez"
$ws.Range("A15").Value = "Ciudad De México"
$ws.Range("A22").Value = "Estado De México"
$ws.Range("A28").Value = "Guanajuato"
$ws.Range("B28").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B33").Value = "Ajuchitlán Del Progreso"
$ws.Range("B34").Value = "Alcozauca De Guerero"
$ws.Range("B36").Value = "Ayutla De Los Libres"
$ws.Range("B37").Value = "Chilapa De Álvarez"
$ws.Range("B38").Value = "Chilpancingo De Los Bravo"
$ws.Range("B39").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B41").Value = "Tlapa De Comonfort"
$ws.Range("B51").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B55").Value = "Oaxaca De Juárez"
$ws.Range("B61").Value = "Villa De Etla"
$ws.Range("B65").Value = "Izúcar De Matamoros"
$ws.Range("B72").Value = "Tepanco De López"
$ws.Range("B74").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B75").Value = "Xayacatlán De Bravo"
$ws.Range("B83").Value = "Apetatitlán De Antonio Carvajal"

# Remove trailing metadata/footnote rows (98-102); also drops the stale
# dimension extent down to the real data range (A1:D96)
$ws.Rows("98:102").Delete()
